$wb = $excel.ActiveWorkbook

# --- Add the new "ageLabel" column (D) to the dataEducation sheet ---
$wsEdu = $wb.Worksheets.Item("dataEducation")
$wsEdu.Range("D1").Value = "ageLabel"
$wsEdu.Range("D2").Value = "25 - 34"
$wsEdu.Range("D3").Value = "35 - 44"
$wsEdu.Range("D4").Value = "45 - 64"
$wsEdu.Range("D5").Value = "65 - 999"

# --- Update selections/active sheet to match the final saved state ---
$wsData = $wb.Worksheets.Item("data")
[void]$wsData.Select()
[void]$wsData.Range("D5:D8").Select()

$wsRef = $wb.Worksheets.Item("references")
[void]$wsRef.Select()
[void]$wsRef.Range("A20").Select()

# dataEducation ends up as the last-active (tabSelected) sheet, cell D7 selected
[void]$wsEdu.Select()
[void]$wsEdu.Range("D7").Select()
